$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 279, shifting existing rows 279-340 down to 280-341.
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with the new weekly price record.
$ws.Range("A279").Value = 10
$ws.Range("B279").Value = "Vega Modelo de Temuco"
$ws.Range("C279").Value = "La Araucanía"
$ws.Range("D279").Value = 45211
$ws.Range("D279").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E279").Value = 9
$ws.Range("F279").Value = 100112005
$ws.Range("G279").Value = "Puerro"
$ws.Range("H279").Value = "Azul de Maquehue"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 100
$ws.Range("K279").Value = 8000
$ws.Range("L279").Value = 8000
$ws.Range("M279").Value = 8000
$ws.Range("N279").Value = "$/docena de paquetes"
$ws.Range("O279").Value = "Provincia de Cautín"
$ws.Range("P279").Value = 667
$ws.Range("Q279").Value = 12
$ws.Range("R279").Value = "Hortaliza"
